$d = $word.ActiveDocument

# --- Change 1: "The Seattle, Washington-based startup does have one big
# advantage. The app neatly taps..." -> "LinkMeUp neatly taps..."
$d.Content.Find.Execute(
    "The Seattle, Washington-based startup does have one big advantage. The app neatly",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "LinkMeUp neatly", 2) | Out-Null

# --- Change 2: remove the existing _GoBack bookmark near the end of the
# document (it will be re-created at the "se|parate apps for" split below).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Change 3: "space that LinkMeUp " -> "space that the Seattle-based startup "
$d.Content.Find.Execute(
    "space that LinkMeUp ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "space that the Seattle-based startup ", 2) | Out-Null

# --- Change 4: split "separate apps for " into "se" + bookmark + "parate apps for "
$rng = $d.Content
$rng.Find.Execute("separate apps for ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $rng.Start + 2
$bm = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bm) | Out-Null
